# Volume Fraction.xlsx — rename sheet, strip "Ti6242-1.1-" prefix from the
# CS-* headers, add ten new "Top-*" columns (K:T) with a bold header style,
# freeze the header row, and resize the columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Sheet name: "Sheet1" -> "Ti6242-1.1"
$ws.Name = "Ti6242-1.1"

# 2. Strip the "Ti6242-1.1-" prefix from the existing CS-1..CS-10 headers
#    (A1:J1). Column order stays the same; only the label text changes.
$csHeaders = @(
    "CS-1-(100x)",
    "CS-2-(100x)",
    "CS-3-(100x)",
    "CS-4-(200x)",
    "CS-5-(200x)",
    "CS-6-(200x)",
    "CS-7-(500x)",
    "CS-8-(500x)",
    "CS-9-(500x)",
    "CS-10-(500x)"
)
for ($i = 0; $i -lt $csHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $csHeaders[$i]
}

# 3. Add the ten new "Top-*" headers in K1:T1.
$topHeaders = @(
    "Top-1-(100x)",
    "Top-2-(100x)",
    "Top-3-(100x)",
    "Top-4-(200x)",
    "Top-5-(200x)",
    "Top-6-(200x)",
    "Top-7-(500x)",
    "Top-8-(500x)",
    "Top-9-(500x)",
    "Top-10-(500x)"
)
for ($i = 0; $i -lt $topHeaders.Length; $i++) {
    $ws.Cells.Item(1, 10 + $i + 1).Value = $topHeaders[$i]
}

# 4. Bold the whole header row (A1:T1) -- new font + cellXf in styles.xml.
$ws.Range("A1:T1").Font.Bold = $true

# 5. Freeze the header row and park the selection on D8, matching the
#    saved view state.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D8").Select()

# 6. Resize the columns:
#      A:I  -> 11
#      J    -> 12
#      K:S  -> 12.140625
#      T    -> 13.28515625
#    ColumnWidth is stored as "characters"; Excel re-expresses it on save
#    as (round(chars*6)/6 + 5/6). Pick chars values that land exactly on
#    the desired stored widths (11 and 12 land exactly; 12.140625 and
#    13.28515625 land on the nearest reachable grid point).
$ws.Columns("A:I").ColumnWidth = 10.166666666666666
$ws.Columns("J:J").ColumnWidth = 11.166666666666666
$ws.Columns("K:S").ColumnWidth = 11.333333333333334
$ws.Columns("T:T").ColumnWidth = 12.5
